$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D" = 0.0214
    "G" = 0.1325301204819277
    "H" = 0.1325301204819277
    "I" = 0.1157356699525374
    "J" = 0.09448867382692233
    "K" = 54.5
    "L" = 0.09948886454910552
    "M" = 0.007
    "N" = 0.00001603665521191294
    "O" = 0.0001284403669724771
    "P" = 0.007
    "Q" = 0.00001603665521191294
    "R" = 0.0001284403669724771
    "U" = 107
    "V" = 0.2451317296678121
    "W" = 0.09903688896965293
    "X" = 0.06179411626139064
    "Y" = 0.03724277270826229
    "Z" = 1.118404262920959
    "AA" = 0.1056765356057779
    "AB" = 0.05842742101559856
    "AC" = 0.04724911459017938
    "AD" = 59.8
    "AF" = 59.8
    "AG" = -47.2
    "AH" = 0.1204916381221036
    "AI" = 0.08755490483162517
    "AJ" = -0.1212432571281788
    "AK" = -0.08194444444444444
    "AL" = 1.84
    "AM" = 1.84
    "AN" = 0.8542857142857142
    "AO" = 34.45652173913043
    "AP" = -0.6742857142857144
    "AQ" = 34.45652173913043
}

foreach ($row in 2..3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
